# Updated model run for baseline
# Apply the diff: update Kernel Size (col C) for rows 2-9 to include an
# extra (3, 3) tuple, and update the Train/Validation Accuracy values
# (cols L, M) for the affected rows to reflect the new model run results.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (Kernel Size) changes for rows 2-9
foreach ($r in 2..9) {
    $ws.Range("C$r").Value = "[(3, 3), (3, 3), (3, 3)]"
}

# Row 2: Validation Accuracy (M2) 1 -> 0.9833333492279053
$ws.Range("M2").Value = 0.9833333492279053

# Row 4: Validation Accuracy (M4) 0.949999988079071 -> 0.9833333492279053
$ws.Range("M4").Value = 0.9833333492279053

# Row 6: Train Accuracy (L6) 1 -> 0.9937499761581421
#        Validation Accuracy (M6) 0.9833333492279053 -> 1
$ws.Range("L6").Value = 0.9937499761581421
$ws.Range("M6").Value = 1

# Row 7: Validation Accuracy (M7) 1 -> 0.949999988079071
$ws.Range("M7").Value = 0.949999988079071

# Row 8: Train Accuracy (L8) 0.9937499761581421 -> 0.9916666746139526
#        Validation Accuracy (M8) 1 -> 0.9333333373069763
$ws.Range("L8").Value = 0.9916666746139526
$ws.Range("M8").Value = 0.9333333373069763
